$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("C").Delete()
$ws.Range("E7").Select() | Out-Null
